# Insert two new paragraphs (a blank line, then a "گام دوم" step heading)
# right after the paragraph that ends with the sentence about the Docker
# image list being verified, and right before the short tab-only paragraph
# that follows it.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive text rather than a fixed
# index, so the script is resilient to any earlier content differences.
$marker = "همانطور که مشاهده میشود در لیست تصاویز داکری"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains($marker)) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the target paragraph (Docker image list sentence)."
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$nextPara = $targetPara.Next()
$insertionPoint = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# First new paragraph: an empty line, formatting matches the paragraph that
# used to directly follow the target (tabs + bidi, no rtl on the mark).
$blankPara = "<w:p $w>" +
    "<w:pPr>" +
        "<w:tabs><w:tab w:val=`"left`" w:pos=`"6147`"/></w:tabs>" +
        "<w:bidi/>" +
        "<w:rPr>" +
            "<w:rFonts w:cs=`"B Nazanin`"/>" +
            "<w:sz w:val=`"28`"/>" +
            "<w:szCs w:val=`"28`"/>" +
            "<w:lang w:bidi=`"fa-IR`"/>" +
        "</w:rPr>" +
    "</w:pPr>" +
"</w:p>"

# Second new paragraph: the bold "گام دوم" (step two) heading, styled like
# the existing "گام اول" heading earlier in the document.
$stepPara = "<w:p $w>" +
    "<w:pPr>" +
        "<w:tabs><w:tab w:val=`"left`" w:pos=`"6147`"/></w:tabs>" +
        "<w:bidi/>" +
        "<w:rPr>" +
            "<w:rFonts w:cs=`"B Nazanin+ Bold`"/>" +
            "<w:sz w:val=`"28`"/>" +
            "<w:szCs w:val=`"28`"/>" +
            "<w:rtl/>" +
            "<w:lang w:bidi=`"fa-IR`"/>" +
        "</w:rPr>" +
    "</w:pPr>" +
    "<w:r>" +
        "<w:rPr>" +
            "<w:rFonts w:cs=`"B Nazanin+ Bold`" w:hint=`"cs`"/>" +
            "<w:sz w:val=`"28`"/>" +
            "<w:szCs w:val=`"28`"/>" +
            "<w:rtl/>" +
            "<w:lang w:bidi=`"fa-IR`"/>" +
        "</w:rPr>" +
        "<w:t>گام دوم</w:t>" +
    "</w:r>" +
"</w:p>"

# A trailing empty paragraph terminator. InsertXML merges the *last*
# paragraph of the fragment into the host paragraph when that last
# paragraph carries run content, so an empty terminator keeps the two
# paragraphs above as fully standalone; the terminator itself is removed
# right afterwards.
$terminator = "<w:p $w></w:p>"

$beforeCount = $d.Paragraphs.Count
$insertionPoint.InsertXML($blankPara + $stepPara + $terminator)

# The terminator landed as its own empty paragraph right after the new
# "گام دوم" paragraph (i.e. right before the original next paragraph) -
# remove it, leaving only the two intended new paragraphs behind.
$strayIndex = $targetIndex + 3
$strayPara = $d.Paragraphs.Item($strayIndex)
$strayPara.Range.Delete()

Write-Output "paragraphs before=$beforeCount after=$($d.Paragraphs.Count)"
